$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 196
$wsOff.Range("C2").Value = 137
$wsOff.Range("D2").Value = 54
$wsOff.Range("E2").Value = 26
$wsOff.Range("F2").Value = 4

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 203
$wsDef.Range("C2").Value = 137
$wsDef.Range("D2").Value = 50
$wsDef.Range("E2").Value = 26
$wsDef.Range("F2").Value = 5
$wsDef.Range("G2").Value = 3
